# "arrow label and camera"
# Remove the blank data row (old row 15, between the "plank" block and
# "scharnier"), which shifts "scharnier"/"slot" rows up and updates several
# quantities/lengths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank row that used to sit at row 15 (between the "plank" group and
# "scharnier"). This shifts scharnier (16->15) and slot (17->16) up by one row,
# and Excel automatically re-adjusts the dependent merged ranges (A7:A15->A7:A14,
# B7:B15->B7:B14, C12:C15->C12:C14).
$ws.Rows.Item(15).Delete()

# The C-column merge that used to cover C10:C11 needs to move up to C9:C10 (the
# "breedte" value moves from the old C10 cell up into C9, and a new value is
# entered in the now-unmerged C11).
$ws.Range("C10:C11").UnMerge()
$ws.Range("C9:C10").Merge()
# Merging clears the lower cell's content (matches target: C10 ends up blank)
# and normalize the border styling back to the plain table style.
$ws.Range("C9:C11").Borders.LineStyle = 1

# --- balk group (rows 2-5) ---
$ws.Range("D2").Value = 39.6
$ws.Range("E2").Value = 4

$ws.Range("D3").Value = 45.6
$ws.Range("E3").Value = 4

$ws.Range("D4").Value = 84.59999999999999

$ws.Range("D5").Value = 89.59999999999999
$ws.Range("E5").Value = 2

# --- blok row (row 6) ---
$ws.Range("E6").Value = 4

# --- plank group (rows 7-14) ---
$ws.Range("D7").Value = 95
$ws.Range("E7").Value = 2

$ws.Range("C8").Value = 12.2
$ws.Range("D8").Value = 95
$ws.Range("E8").Value = 2

$ws.Range("C9").Value = 12.8
$ws.Range("D9").Value = 95

$ws.Range("D10").Value = 95.59999999999999

$ws.Range("C11").Value = 19.8
$ws.Range("D11").Value = 94.59999999999999

$ws.Range("D12").Value = 76
$ws.Range("E12").Value = 3

$ws.Range("D13").Value = 95
$ws.Range("E13").Value = 10

$ws.Range("D14").Value = 95.59999999999999
$ws.Range("E14").Value = 2

# --- scharnier row (now row 15) ---
$ws.Range("E15").Value = 3

# --- slot row (now row 16) ---
$ws.Range("E16").Value = 1
